$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.395
$ws.Range("G2").Value = -0.108843537414966
$ws.Range("H2").Value = -0.108843537414966
$ws.Range("I2").Value = -0.2057823129251701
$ws.Range("J2").Value = -0.1617364411518032
$ws.Range("K2").Value = 1.221
$ws.Range("L2").Value = 2.076530612244898
$ws.Range("O2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("U2").Value = 5.067
$ws.Range("V2").Value = 0.01441126279863481
$ws.Range("W2").Value = 0.1913035381750466
$ws.Range("X2").Value = 0.04537944222634582
$ws.Range("Y2").Value = 0.1459240959487008
$ws.Range("Z2").Value = 0.01175811870100784
$ws.Range("AA2").Value = 0.1655087971465836
$ws.Range("AB2").Value = 0.04537944222634582
$ws.Range("AC2").Value = 0.1174911671881717
$ws.Range("AD2").Value = 34.5
$ws.Range("AF2").Value = 34.5
$ws.Range("AG2").Value = 29.433
$ws.Range("AH2").Value = 0.08935508935508936
$ws.Range("AI2").Value = 0.5706912808297354
$ws.Range("AJ2").Value = 0.07724527796805002
$ws.Range("AK2").Value = 0.531415881269635
$ws.Range("AM2").Value = -0.045
$ws.Range("AN2").Value = -453.9473684210527
$ws.Range("AP2").Value = -387.2763157894737
$ws.Range("AQ2").Value = 2.688888888888889

# Row 3
$ws.Range("D3").Value = 0.395
$ws.Range("I3").Value = 0.1666666666666667
$ws.Range("J3").Value = 0.1073059360730594
$ws.Range("K3").Value = 0.047
$ws.Range("L3").Value = 0.07993197278911565
$ws.Range("U3").Value = 0.153
$ws.Range("V3").Value = 0.0006139646869983949
$ws.Range("W3").Value = 0.3133333333333334
$ws.Range("X3").Value = 0.04537944222634582
$ws.Range("Y3").Value = 0.2679538911069875
$ws.Range("Z3").Value = 117.5999999999999
$ws.Range("AA3").Value = 12.61917808219177
$ws.Range("AB3").Value = 0.04537944222634582
$ws.Range("AC3").Value = 12.57379863996542
$ws.Range("AG3").Value = -0.153
$ws.Range("AJ3").Value = -0.0006143418712130642
$ws.Range("AK3").Value = -2.833333333333333

# Row 4
$ws.Range("K4").Value = -0.063
$ws.Range("U4").Value = 0.027
$ws.Range("V4").Value = 0.00125
$ws.Range("W4").Value = 0.3641618497109827
$ws.Range("X4").Value = 0.04537944222634582
$ws.Range("Y4").Value = 0.3187824074846369
$ws.Range("AA4").Value = 0.3333333333333333
$ws.Range("AB4").Value = 0.04537944222634582
$ws.Range("AC4").Value = 0.2879538911069875
$ws.Range("AG4").Value = -0.027
$ws.Range("AI4").Value = 0
$ws.Range("AJ4").Value = -0.001251564455569462
$ws.Range("AK4").Value = -0.07317073170731707

# Row 5
$ws.Range("K5").Value = 1.24
$ws.Range("O5").Value = 0
$ws.Range("R5").Value = 0
$ws.Range("U5").Value = 0.017
$ws.Range("V5").Value = 0.0005862068965517242
$ws.Range("W5").Value = 0.06927374301675979
$ws.Range("X5").Value = 0.07542490877000968
$ws.Range("Y5").Value = -0.006151165753249893
$ws.Range("AA5").Value = -0.002315739040166095
$ws.Range("AB5").Value = 0.05065581769047794
$ws.Range("AC5").Value = -0.05297155673064404
$ws.Range("AD5").Value = 34.5
$ws.Range("AF5").Value = 34.5
$ws.Range("AG5").Value = 34.483
$ws.Range("AH5").Value = 0.5433070866141733
$ws.Range("AI5").Value = 0.6284153005464481
$ws.Range("AJ5").Value = 0.5431847896287195
$ws.Range("AK5").Value = 0.6283002022484193
$ws.Range("AN5").Value = -453.9473684210527
$ws.Range("AP5").Value = -453.7236842105263

# Row 6
$ws.Range("K6").Value = -0.003
$ws.Range("U6").Value = 4.87
$ws.Range("V6").Value = 0.09401544401544402
$ws.Range("W6").Value = -0.0006479481641468683
$ws.Range("X6").Value = 0.04537944222634582
$ws.Range("Y6").Value = -0.04602739039049269
$ws.Range("AA6").Value = -0.2000000000000007
$ws.Range("AB6").Value = 0.04537944222634582
$ws.Range("AC6").Value = -0.2453794422263465
$ws.Range("AG6").Value = -4.87
$ws.Range("AJ6").Value = -0.1037715746857021
$ws.Range("AK6").Value = -60.87499999999995
$ws.Range("AM6").Value = -0.045
$ws.Range("AQ6").Value = 0.888888888888889
